$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 (HOLA / ADIOS / QUE TAL ESTAS? / A MI ...)
$ws.Range("A5").Value = "HOLA"
$ws.Range("B5").Value = "ADIOS"
$ws.Range("C5").Value = "QUE TAL ESTAS?"
$ws.Range("D5").Value = "A MI"

# New column E (SABER / CONTAR / VIVIR / SENTIR / DECIR)
$ws.Range("E1").Value = "SABER"
$ws.Range("E2").Value = "CONTAR"
$ws.Range("E3").Value = "VIVIR"
$ws.Range("E4").Value = "SENTIR"
$ws.Range("E5").Value = "DECIR"

# Update the active selection to match the saved view
$ws.Range("F5").Select() | Out-Null
